# Auto-generated Excel COM-interop script to apply numeric updates
# to the Typhon_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each affected row, set the new values for changed cells, remove cells
# that were dropped entirely (ClearContents), and populate newly added cells.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 378.22223
$ws.Range("I98").Value = 378
$ws.Range("K98").Value = 378
$ws.Range("M98").Value = 1120
# Row 122
$ws.Range("H122").Value = 378.22223
$ws.Range("I122").Value = 378
$ws.Range("K122").Value = 1134
$ws.Range("M122").Value = 1316
# Row 129
$ws.Range("H129").Value = 676.5714
# Row 132
$ws.Range("H132").Value = 44750.957
$ws.Range("I132").Value = 44750.957
$ws.Range("K132").Value = 134252.871
$ws.Range("M132").Value = -131722.871
# Row 138
$ws.Range("H138").Value = 2504.7073
$ws.Range("J138").Value = 2958.3635
$ws.Range("L138").Value = 8875.0905
$ws.Range("N138").Value = -19155.0905

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 28716.074
$ws.Range("I32").Value = 37014.793
$ws.Range("J32").Value = 6837.636
$ws.Range("K32").Value = 37014.793
$ws.Range("L32").Value = 6837.636
$ws.Range("M32").Value = -36727.793
$ws.Range("N32").Value = -7411.636
# Row 45
$ws.Range("H45").Value = 2666.6191
$ws.Range("I45").Value = 2168.7
$ws.Range("J45").Value = 3119.2727
$ws.Range("K45").Value = 2168.7
$ws.Range("L45").Value = 3119.2727
$ws.Range("M45").Value = -1791.7
$ws.Range("N45").Value = -3873.2727
# Row 59
$ws.Range("H59").Value = 22000
$ws.Range("J59").Value = 22000
$ws.Range("L59").Value = 22000
$ws.Range("N59").Value = -23608
# Row 74
$ws.Range("H74").Value = 958.6
$ws.Range("I74").Value = 511.0909
$ws.Range("K74").Value = 511.0909
$ws.Range("M74").Value = 362.9091
# Row 77
$ws.Range("H77").Value = 958.6
$ws.Range("I77").Value = 511.0909
$ws.Range("K77").Value = 2555.4545
$ws.Range("M77").Value = 1812.5455
# Row 122
$ws.Range("H122").Value = 2047.375
$ws.Range("I122").Value = 1783.7142
$ws.Range("K122").Value = 5351.142599999999
$ws.Range("M122").Value = -2901.142599999999
# Row 132
$ws.Range("H132").Value = 23680.912
$ws.Range("I132").Value = 1526.9333
$ws.Range("J132").Value = 65219.625
$ws.Range("K132").Value = 4580.7999
$ws.Range("L132").Value = 195658.875
$ws.Range("M132").Value = -2050.7999
$ws.Range("N132").Value = -200718.875

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 29366.945
$ws.Range("I86").Value = 42883.832
$ws.Range("K86").Value = 42883.832
$ws.Range("M86").Value = -41760.832
# Row 89
$ws.Range("H89").Value = 29366.945
$ws.Range("I89").Value = 42883.832
$ws.Range("K89").Value = 214419.16
$ws.Range("M89").Value = -208803.16
# Row 105
$ws.Range("H105").Value = 1858.7646
$ws.Range("I105").Value = 1599.875
$ws.Range("K105").Value = 1599.875
$ws.Range("M105").Value = 147.125
# Row 107
$ws.Range("H107").Value = 1244.4286
$ws.Range("I107").Value = 1244.4286
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1244.4286
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 675.5714
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9200.071
$ws.Range("I31").Value = 12315.27
$ws.Range("J31").Value = 4137.875
$ws.Range("K31").Value = 12315.27
$ws.Range("L31").Value = 4137.875
$ws.Range("N31").Value = -4727.875
$ws.Range("M31").Value = -12020.27
# Row 34
$ws.Range("H34").Value = 9200.071
$ws.Range("I34").Value = 12315.27
$ws.Range("J34").Value = 4137.875
$ws.Range("K34").Value = 12315.27
$ws.Range("L34").Value = 4137.875
$ws.Range("N34").Value = -4541.875
$ws.Range("M34").Value = -12113.27
# Row 94
$ws.Range("H94").Value = 6347.4375
$ws.Range("J94").Value = 8685.9
$ws.Range("L94").Value = 8685.9
$ws.Range("N94").Value = -9587.9
# Row 122
$ws.Range("H122").Value = 898.2308
$ws.Range("I122").Value = 1088.7778
$ws.Range("J122").Value = 469.5
$ws.Range("K122").Value = 3266.3334
$ws.Range("L122").Value = 1408.5
$ws.Range("M122").Value = -816.3334000000004
$ws.Range("N122").Value = -6308.5
# Row 132
$ws.Range("H132").Value = 55883.1
$ws.Range("I132").Value = 250900.5
$ws.Range("J132").Value = 7128.75
$ws.Range("K132").Value = 752701.5
$ws.Range("L132").Value = 21386.25
$ws.Range("M132").Value = -750171.5
$ws.Range("N132").Value = -26446.25

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1172.9584
$ws.Range("I68").Value = 475.25
$ws.Range("J68").Value = 1312.5
$ws.Range("K68").Value = 1425.75
$ws.Range("L68").Value = 3937.5
$ws.Range("M68").Value = -614.75
$ws.Range("N68").Value = -5559.5
# Row 71
$ws.Range("H71").Value = 1172.9584
$ws.Range("I71").Value = 475.25
$ws.Range("J71").Value = 1312.5
$ws.Range("K71").Value = 4277.25
$ws.Range("L71").Value = 11812.5
$ws.Range("M71").Value = -221.25
$ws.Range("N71").Value = -19924.5
# Row 109
$ws.Range("H109").Value = 5836.75
$ws.Range("J109").Value = 6284.857
$ws.Range("L109").Value = 18854.571
$ws.Range("N109").Value = -20934.571
# Row 131
$ws.Range("H131").Value = 115775.484
$ws.Range("I131").Value = 819
$ws.Range("J131").Value = 125834.17
$ws.Range("K131").Value = 2457
$ws.Range("L131").Value = 377502.51
$ws.Range("M131").Value = 2583
$ws.Range("N131").Value = -387582.51

$ws = $wb.Worksheets.Item("GSM")
# Row 53
$ws.Range("H53").Value = 5039
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
# Row 70
$ws.Range("H70").Value = 4269.3
$ws.Range("J70").Value = 4474.8335
$ws.Range("L70").Value = 4474.8335
$ws.Range("N70").Value = -5014.8335
# Row 73
$ws.Range("H73").Value = 4269.3
$ws.Range("J73").Value = 4474.8335
$ws.Range("L73").Value = 4474.8335
$ws.Range("N73").Value = -6346.8335
# Row 97
$ws.Range("H97").Value = 2020.6364
$ws.Range("I97").Value = 1108.8667
$ws.Range("J97").Value = 3974.4285
$ws.Range("K97").Value = 1108.8667
$ws.Range("L97").Value = 3974.4285
$ws.Range("M97").Value = -612.8667
$ws.Range("N97").Value = -4966.4285
# Row 126
$ws.Range("H126").Value = 4809.304
$ws.Range("I126").Value = 4008.3333
$ws.Range("J126").Value = 5683.091
$ws.Range("K126").Value = 12024.9999
$ws.Range("L126").Value = 17049.273
$ws.Range("M126").Value = -9554.999899999999
$ws.Range("N126").Value = -21989.273
# Row 132
$ws.Range("H132").Value = 70257.95
$ws.Range("I132").Value = 68138.47
$ws.Range("J132").Value = 74799.71000000001
$ws.Range("K132").Value = 204415.41
$ws.Range("L132").Value = 224399.13
$ws.Range("M132").Value = -201885.41
$ws.Range("N132").Value = -229459.13

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1115.4546
$ws.Range("I16").Value = 696.6667
$ws.Range("K16").Value = 696.6667
$ws.Range("M16").Value = -526.6667
# Row 57
$ws.Range("H57").Value = 980
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 4545.4546
$ws.Range("J96").Value = 5055.5557
$ws.Range("L96").Value = 5055.5557
$ws.Range("N96").Value = -7801.5557
# Row 107
$ws.Range("H107").Value = 1389.238
$ws.Range("I107").Value = 802.4
$ws.Range("J107").Value = 1922.7273
$ws.Range("K107").Value = 2407.2
$ws.Range("L107").Value = 5768.1819
$ws.Range("M107").Value = -487.1999999999998
$ws.Range("N107").Value = -9608.1819
# Row 132
$ws.Range("H132").Value = 2539.1765
$ws.Range("I132").Value = 2330.8
$ws.Range("J132").Value = 2836.8572
$ws.Range("K132").Value = 6992.400000000001
$ws.Range("L132").Value = 8510.571599999999
$ws.Range("M132").Value = -4462.400000000001
$ws.Range("N132").Value = -13570.5716
# Row 136
$ws.Range("H136").Value = 2689971.8
$ws.Range("I136").Value = 8065141
$ws.Range("K136").Value = 24195423
$ws.Range("M136").Value = -24192873
